# CV Power BI template
# Adds an "order" ranking column (used to sort entries in the Power BI
# report) to every data sheet, fills in the missing employer for the
# self-employment row on the English "experience" sheet, and widens the
# new column on "education". Finishes by restoring the same active
# sheet/selection layout recorded in the target workbook.

$wb = $excel.ActiveWorkbook

# --- education ---------------------------------------------------------
$ws = $wb.Worksheets.Item("education")
$ws.Range("D1").Value = "order"
$ws.Range("D2").Value = 3
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 1
$ws.Columns.Item(3).ColumnWidth = 9.6

# --- experience ----------------------------------------------------------
$ws = $wb.Worksheets.Item("experience")
$ws.Range("E1").Value = "order"
$ws.Range("E2").Value = 3
$ws.Range("E3").Value = 2
$ws.Range("B4").Value = "Self-emplyment"
$ws.Range("E4").Value = 1

# --- wykształcenie -------------------------------------------------------
$ws = $wb.Worksheets.Item("wykształcenie")
$ws.Range("D1").Value = "order"
$ws.Range("D2").Value = 3
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 1

# --- doświadczenie ---------------------------------------------------------
$ws = $wb.Worksheets.Item("doświadczenie")
$ws.Range("E1").Value = "order"
$ws.Range("E2").Value = 3
$ws.Range("E3").Value = 2
$ws.Range("E4").Value = 1

# --- restore per-sheet selections, leaving "education" as the active tab ---
$ws = $wb.Worksheets.Item("experience")
$ws.Activate()
$ws.Range("C7").Select()

$ws = $wb.Worksheets.Item("wykształcenie")
$ws.Activate()
$ws.Range("C8").Select()

$ws = $wb.Worksheets.Item("doświadczenie")
$ws.Activate()
$ws.Range("E8").Select()

$ws = $wb.Worksheets.Item("education")
$ws.Activate()
$ws.Range("D8").Select()
